$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.505.31'
$ws.Range('E2').Value = '  -1.76%  '
$ws.Range('D3').Value = '1.748.68'
$ws.Range('E3').Value = '  -1.79%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.75'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4466'
$ws.Range('E7').Value = '  +4.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3599'
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07495'
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.99'
$ws.Range('E10').Value = '  -5.53%  '
$ws.Range('E11').Value = '  -1.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.68'
$ws.Range('E13').Value = '  -4.39%  '
$ws.Range('E14').Value = '  -2.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.124'
$ws.Range('E15').Value = '  -2.95%  '
$ws.Range('D16').Value = '1.748.20'
$ws.Range('E16').Value = '  -2.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.65'
$ws.Range('E17').Value = '  +2.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001062'
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06384'
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('E21').Value = '  -2.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.848'
$ws.Range('E22').Value = '  -1.71%  '
$ws.Range('D23').Value = '27.549.35'
$ws.Range('E23').Value = '  -1.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.23'
$ws.Range('E24').Value = '  -1.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.086'
$ws.Range('E25').Value = '  -2.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.27'
$ws.Range('E26').Value = '  +1.23%  '
$ws.Range('E27').Value = '  +0.26%  '
$ws.Range('D28').Value = '1.949.77'
$ws.Range('E29').Value = '  -4.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.66'
$ws.Range('E30').Value = '  -0.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.083'
$ws.Range('E31').Value = '  -7.22%  '
$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.671'
$ws.Range('E32').Value = '  +5.02%  '
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09073'
$ws.Range('E33').Value = '  +1.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.540'
$ws.Range('E34').Value = '  -2.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '11.97'
$ws.Range('E35').Value = '  -5.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02294'
$ws.Range('E36').Value = '  -1.71%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06017'
$ws.Range('E37').Value = '  -0.83%  '
$ws.Range('B38').Value = 'TheSandbox'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6366'
$ws.Range('E38').Value = '  -1.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2087'
$ws.Range('E39').Value = '  -1.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.936'
$ws.Range('E40').Value = '  -2.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.202'
$ws.Range('E41').Value = '  +1.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.383'
$ws.Range('E42').Value = '  -1.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.768'
$ws.Range('E43').Value = '  -1.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.23'
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5899'
$ws.Range('E46').Value = '  -1.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '122.51'
$ws.Range('E47').Value = '  -1.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.955'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.147'
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06856'
$ws.Range('E50').Value = '  -1.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.24'
$ws.Range('E51').Value = '  -3.28%  '
